$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "dfsfs"
$ws.Range("A2").Value = "sdfs"
$ws.Range("A3").Value = "sdfsdf"
$ws.Range("A4").Value = "sfsff"

$ws.Range("A4").Select()
